$wb = $excel.ActiveWorkbook

# ==================== Sheet 1: LP1912 ====================
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: 05:24:16"
$ws1.Range("A3").Value = "Total filas: 29"

$data1 = @(
    @("05:24:16", "05:25", "23_HERNANDEZ", 1, "LP1912"),
    @("04:56:30", "05:34", "215B_EL PATO", 38, "LP1912"),
    @("05:24:16", "05:35", "215B_EL PATO", 11, "LP1912"),
    @("05:24:16", "05:46", "15_ABASTO", 22, "LP1912"),
    @("05:24:16", "05:54", "10_OLMOS", 30, "LP1912"),
    @("05:24:16", "06:04", "16_SANTA ANA", 40, "LP1912"),
    @("05:24:16", "06:11", "215A_EL PATO", 47, "LP1912"),
    @("05:24:16", "06:14", "225_HARAS DEL SUR", 50, "LP1912"),
    @("04:56:30", "06:18", "16_SANTA ANA", 82, "LP1912"),
    @("05:24:16", "06:21", "26_HERNANDEZ", 57, "LP1912"),
    @("04:45:48", "06:24", "16_SANTA ANA", 99, "LP1912"),
    @("05:24:16", "06:27", "23_HERNANDEZ", 63, "LP1912"),
    @("04:56:30", "06:29", "86_EST CHICA-ESC AGRARIA", 93, "LP1912"),
    @("05:24:16", "06:30", "86_EST CHICA-ESC AGRARIA", 66, "LP1912"),
    @("05:24:16", "06:31", "16_SANTA ANA", 67, "LP1912"),
    @("05:24:16", "06:44", "225_C ROCA-H SUR", 80, "LP1912"),
    @("05:24:16", "06:46", "215C_EL PATO", 82, "LP1912"),
    @("05:24:16", "07:00", "14_ABASTO", 96, "LP1912"),
    @("05:24:16", "07:05", "15_ABASTO", 101, "LP1912"),
    @("05:24:16", "07:07", "225_GOMEZ", 103, "LP1912"),
    @("05:24:16", "07:11", "215A_EL PATO", 107, "LP1912"),
    @("05:24:16", "07:16", "11_ETCHEVERRY", 112, "LP1912"),
    @("05:24:16", "07:21", "26_HERNANDEZ", 117, "LP1912"),
    @("05:24:16", "07:23", "10_OLMOS", 119, "LP1912")
)
$r = 11
foreach ($row in $data1) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# ==================== Sheet 2: LP1912-215 ====================
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 05:24:16"
$ws2.Range("A3").Value = "Total filas: 7"

$ws2.Cells.Item(9, 1).Value = "05:24:16"
$ws2.Cells.Item(9, 2).Value = "05:35"
$ws2.Cells.Item(9, 3).Value = "215B_EL PATO"
$ws2.Cells.Item(9, 4).Value = 11
$ws2.Cells.Item(9, 5).Value = "LP1912"

$ws2.Cells.Item(10, 1).Value = "05:24:16"
$ws2.Cells.Item(10, 2).Value = "06:11"
$ws2.Cells.Item(10, 3).Value = "215A_EL PATO"
$ws2.Cells.Item(10, 4).Value = 47
$ws2.Cells.Item(10, 5).Value = "LP1912"

$ws2.Cells.Item(11, 1).Value = "05:24:16"
$ws2.Cells.Item(11, 2).Value = "06:46"
$ws2.Cells.Item(11, 3).Value = "215C_EL PATO"
$ws2.Cells.Item(11, 4).Value = 82
$ws2.Cells.Item(11, 5).Value = "LP1912"

$ws2.Cells.Item(12, 1).Value = "05:24:16"
$ws2.Cells.Item(12, 2).Value = "07:11"
$ws2.Cells.Item(12, 3).Value = "215A_EL PATO"
$ws2.Cells.Item(12, 4).Value = 107
$ws2.Cells.Item(12, 5).Value = "LP1912"

# ==================== Sheet 3: 6203-6173 ====================
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 05:24:16"
$ws3.Range("A3").Value = "Total filas: 7"

$ws3.Cells.Item(7, 1).Value = "05:24:16"
$ws3.Cells.Item(7, 2).Value = "05:44"
$ws3.Cells.Item(7, 3).Value = "215A_LA PLATA"
$ws3.Cells.Item(7, 4).Value = 20
$ws3.Cells.Item(7, 5).Value = "L6173"

$ws3.Cells.Item(9, 1).Value = "05:24:16"
$ws3.Cells.Item(9, 2).Value = "06:09"
$ws3.Cells.Item(9, 3).Value = "215A_LA PLATA"
$ws3.Cells.Item(9, 4).Value = 45
$ws3.Cells.Item(9, 5).Value = "L6173"

$ws3.Cells.Item(11, 1).Value = "05:24:16"
$ws3.Cells.Item(11, 2).Value = "06:33"
$ws3.Cells.Item(11, 3).Value = "215C_LA PLATA"
$ws3.Cells.Item(11, 4).Value = 69
$ws3.Cells.Item(11, 5).Value = "L6203"

$ws3.Cells.Item(12, 1).Value = "05:24:16"
$ws3.Cells.Item(12, 2).Value = "07:00"
$ws3.Cells.Item(12, 3).Value = "215B_LP-P MOR-1 Y 57"
$ws3.Cells.Item(12, 4).Value = 96
$ws3.Cells.Item(12, 5).Value = "L6173"

